$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values on specific rows to reflect repulled/recalculated data
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -2
$ws.Range("F12").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = -13
$ws.Range("F21").Value = -1
$ws.Range("F24").Value = 3
$ws.Range("F33").Value = -6
